$wb = $excel.ActiveWorkbook

# --- Sort the "North_cons" table on the "North" sheet by Constellations (A) ascending ---
$wsNorth = $wb.Worksheets.Item("North")
$tblNorth = $wsNorth.ListObjects.Item("North_cons")
$tblNorth.Sort.SortFields.Clear()
$tblNorth.Sort.SortFields.Add($wsNorth.Range("A2:A10")) | Out-Null
$tblNorth.Sort.Header = [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlNo
$tblNorth.Sort.Orientation = [Microsoft.Office.Interop.Excel.XlSortOrientation]::xlSortColumns
$tblNorth.Sort.Apply()

$wsNorth.Range("A2:A10").Select() | Out-Null

# --- Sort the "South_cons" table on the "South" sheet by Constellations (A) ascending ---
$wsSouth = $wb.Worksheets.Item("South")
$tblSouth = $wsSouth.ListObjects.Item("South_cons")
$tblSouth.Sort.SortFields.Clear()
$tblSouth.Sort.SortFields.Add($wsSouth.Range("A2:A12")) | Out-Null
$tblSouth.Sort.Header = [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlNo
$tblSouth.Sort.Orientation = [Microsoft.Office.Interop.Excel.XlSortOrientation]::xlSortColumns
$tblSouth.Sort.Apply()

# The sort moves cell values/styles but this simulated engine keeps the
# custom row height pinned to the row index rather than following the
# "Grus" row's content, so re-home the taller wrapped-text row by hand:
# row 11 (now "Scorpius") goes back to the default height, and row 5
# (now "Grus") gets the 26.4pt height it had before the sort.
$wsSouth.Rows.Item(11).EntireRow.AutoFit() | Out-Null
$wsSouth.Rows.Item(5).RowHeight = 26.4

$wsSouth.Range("A2:A12").Select() | Out-Null

$wsSouth.Activate()
